$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.424501
$ws.Range("H2").Value = 10.273503
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.616180666666666
$ws.Range("N2").Value = 22.848542
$ws.Range("O2").Value = 0.06266940401417194
$ws.Range("P2").Value = 0.06266940401417194
$ws.Range("Q2").Value = 26.08161830918066
$ws.Range("R2").Value = 234.734564782626
$ws.Range("S2").Value = 0.06266940401417194
$ws.Range("T2").Value = 0.06266940401417194

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.424501
$ws.Range("H3").Value = 10.273503
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 86.98680366666667
$ws.Range("N3").Value = 260.960411
$ws.Range("O3").Value = 0.7157670466966058
$ws.Range("P3").Value = 0.7157670466966058
$ws.Range("Q3").Value = 297.8863961433037
$ws.Range("R3").Value = 2680.977565289733
$ws.Range("S3").Value = 0.7157670466966058
$ws.Range("T3").Value = 0.7157670466966058

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.424501
$ws.Range("H4").Value = 10.273503
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.9360126666666666
$ws.Range("N4").Value = 2.808038
$ws.Range("O4").Value = 0.007701938614251506
$ws.Range("P4").Value = 0.007701938614251506
$ws.Range("Q4").Value = 3.205376313012666
$ws.Range("R4").Value = 28.848386817114
$ws.Range("S4").Value = 0.007701938614251506
$ws.Range("T4").Value = 0.007701938614251506

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.424501
$ws.Range("H5").Value = 10.273503
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.99049233333334
$ws.Range("N5").Value = 77.97147700000001
$ws.Range("O5").Value = 0.2138616106749707
$ws.Range("P5").Value = 0.2138616106749707
$ws.Range("Q5").Value = 89.00446698599234
$ws.Range("R5").Value = 801.0402028739311
$ws.Range("S5").Value = 0.2138616106749707
$ws.Range("T5").Value = 0.2138616106749707
